$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Coin name) updates
$bChanges = @{
    35 = "Stellar"
    36 = "WEMIXTOKEN"
}
foreach ($row in $bChanges.Keys) {
    $ws.Cells.Item($row, 2).Value = $bChanges[$row]
}

# Column C (Link) updates
$cChanges = @{
    35 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    36 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
}
foreach ($row in $cChanges.Keys) {
    $ws.Cells.Item($row, 3).Value = $cChanges[$row]
}

# Column D (Price) updates - text-like values (contain multiple dots, safe as text)
$dTextChanges = @{
    2 = "24.873.05"
    3 = "1.708.30"
    16 = "1.708.09"
    24 = "24.898.88"
    30 = "1.895.76"
}
foreach ($row in $dTextChanges.Keys) {
    $ws.Cells.Item($row, 4).Value = $dTextChanges[$row]
}

# Column D (Price) updates - number-like values (must force text storage to avoid numeric coercion)
$dNumChanges = @{
    5 = "312.25"
    6 = "0.9991"
    7 = "0.3738"
    8 = "49.44"
    9 = "0.3437"
    10 = "1.223"
    11 = "0.07549"
    13 = "21.31"
    14 = "6.326"
    15 = "7.107"
    17 = "0.00001132"
    18 = "0.06734"
    19 = "0.9987"
    20 = "84.09"
    21 = "17.33"
    22 = "6.390"
    23 = "13.12"
    25 = "2.452"
    26 = "2.807"
    27 = "20.41"
    28 = "149.80"
    29 = "133.16"
    31 = "1.246"
    32 = "6.829"
    33 = "4.228"
    34 = "13.90"
    35 = "0.08802"
    36 = "1.769"
    37 = "5.625"
    38 = "0.06659"
    39 = "9.195"
    40 = "0.02416"
    41 = "0.2239"
    42 = "1.278"
    43 = "0.6470"
    44 = "0.9988"
    45 = "13.95"
    46 = "0.6166"
    47 = "3.834"
    48 = "2.132"
    49 = "129.86"
    50 = "0.07321"
    51 = "80.21"
}
foreach ($row in $dNumChanges.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $dNumChanges[$row]
    $cell.Style = "Normal"
}

# Column E (Volume 1h) updates
$eChanges = @{
    2 = "  +1.63%  "
    3 = "  +1.55%  "
    4 = "  +0.07%  "
    5 = "  +1.38%  "
    6 = "  +0.07%  "
    7 = "  +0.63%  "
    8 = "  +3.42%  "
    9 = "  -0.31%  "
    10 = "  +3.78%  "
    11 = "  +3.64%  "
    12 = "  +0.10%  "
    13 = "  +4.64%  "
    14 = "  +3.16%  "
    15 = "  +5.05%  "
    16 = "  +1.60%  "
    17 = "  +2.03%  "
    18 = "  +0.77%  "
    19 = "  +0.00%  "
    21 = "  +4.97%  "
    22 = "  +4.21%  "
    23 = "  +7.33%  "
    24 = "  +1.93%  "
    25 = "  +0.12%  "
    26 = "  +5.08%  "
    27 = "  +4.40%  "
    28 = "  -2.65%  "
    29 = "  +4.54%  "
    30 = "  +1.53%  "
    31 = "  +27.06%  "
    32 = "  +7.73%  "
    33 = "  +4.13%  "
    34 = "  +11.90%  "
    35 = "  +3.70%  "
    36 = "  +3.61%  "
    37 = "  +4.49%  "
    38 = "  +2.32%  "
    39 = "  +3.30%  "
    40 = "  +3.40%  "
    41 = "  +5.74%  "
    42 = "  +1.67%  "
    43 = "  +4.77%  "
    44 = "  +0.08%  "
    45 = "  +5.48%  "
    46 = "  +3.47%  "
    47 = "  +1.45%  "
    48 = "  +5.06%  "
    49 = "  +1.69%  "
    50 = "  +1.76%  "
    51 = "  +5.36%  "
}
foreach ($row in $eChanges.Keys) {
    $ws.Cells.Item($row, 5).Value = $eChanges[$row]
}